# Updates the cryptos worksheet with the latest scraped figures.
# Numeric-looking strings (prices such as "225.64") must stay stored as
# plain text, exactly like the original inline strings, so we briefly
# force a text number-format while assigning the value and then restore
# the cell's original ("Normal") style to avoid leaving stray formatting
# behind.
function Set-TextValue {
    param($ws, [string]$addr, [string]$val)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
Set-TextValue $ws "D2" "34.502.18"
Set-TextValue $ws "E2" "  +0.49%  "

# Row 3 - Ethereum
Set-TextValue $ws "D3" "1.811.43"
Set-TextValue $ws "E3" "  +0.43%  "

# Row 4 - TetherUSD
Set-TextValue $ws "E4" "  -0.15%  "

# Row 5 - BNB
Set-TextValue $ws "D5" "225.64"
Set-TextValue $ws "E5" "  -0.99%  "

# Row 6 - XRP
Set-TextValue $ws "E6" "  +2.84%  "

# Row 7 - USDC
Set-TextValue $ws "E7" "  -0.15%  "

# Row 8 - Solana
Set-TextValue $ws "D8" "38.25"
Set-TextValue $ws "E8" "  +6.17%  "

# Row 9 - Cardano
Set-TextValue $ws "E9" "  -4.00%  "

# Row 10 - Dogecoin
Set-TextValue $ws "E10" "  -2.36%  "

# Row 11 - TRON
Set-TextValue $ws "E11" "  +0.82%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue $ws "D12" "2.072.95"
Set-TextValue $ws "E12" "  +0.46%  "

# Row 13 - Chainlink
Set-TextValue $ws "D13" "11.24"
Set-TextValue $ws "E13" "  -3.52%  "

# Row 14 - WrappedEther
Set-TextValue $ws "D14" "1.807.92"
Set-TextValue $ws "E14" "  +0.47%  "

# Row 15 - Polygon
Set-TextValue $ws "D15" "0.634"
Set-TextValue $ws "E15" "  -1.67%  "

# Row 16 - WrappedBTC
Set-TextValue $ws "D16" "34.469.21"
Set-TextValue $ws "E16" "  +0.38%  "

# Row 17 - Polkadot
Set-TextValue $ws "E17" "  -1.41%  "

# Row 18 - Litecoin
Set-TextValue $ws "D18" "68.37"
Set-TextValue $ws "E18" "  -1.05%  "

# Row 19 - BitcoinCash
Set-TextValue $ws "D19" "243.31"
Set-TextValue $ws "E19" "  -0.88%  "

# Row 20 - ShibaInu
Set-TextValue $ws "E20" "  -2.46%  "

# Row 21 - Avalanche
Set-TextValue $ws "E21" "  -2.47%  "

# Row 22 - Dai
Set-TextValue $ws "E22" "  -0.13%  "

# Row 23 - Uniswap
Set-TextValue $ws "E23" "  -1.32%  "

# Row 24 - Toncoin
Set-TextValue $ws "E24" "  +3.94%  "

# Row 25 - Monero
Set-TextValue $ws "D25" "170.76"
Set-TextValue $ws "E25" "  -1.18%  "

# Row 26 - Cosmos
Set-TextValue $ws "E26" "  -2.21%  "

# Row 27 - EthereumClassic
Set-TextValue $ws "D27" "17.68"
Set-TextValue $ws "E27" "  +4.75%  "

# Row 28 - Stellar
Set-TextValue $ws "E28" "  +2.00%  "

# Row 29 - BinanceUSD
Set-TextValue $ws "E29" "  -0.19%  "

# Row 30 - Filecoin
Set-TextValue $ws "E30" "  -0.96%  "

# Row 31 - PancakeSwap
Set-TextValue $ws "E31" "  -1.53%  "

# Row 32 - Hedera
Set-TextValue $ws "E32" "  -2.74%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue $ws "E33" "  -4.35%  "

# Row 34 - LidoDAOToken
Set-TextValue $ws "D34" "1.82"
Set-TextValue $ws "E34" "  +0.14%  "

# Row 35 - Maker
Set-TextValue $ws "D35" "1.356.87"
Set-TextValue $ws "E35" "  -2.62%  "

# Row 36 - ImmutableX
Set-TextValue $ws "D36" "0.644"
Set-TextValue $ws "E36" "  -4.29%  "

# Row 37 - TrustWalletToken
Set-TextValue $ws "E37" "  -0.89%  "

# Row 38 - VeChain
Set-TextValue $ws "D38" "0.0187"
Set-TextValue $ws "E38" "  -1.51%  "

# Row 39 - RenderToken
Set-TextValue $ws "E39" "  -4.89%  "

# Row 40 - HuobiToken
Set-TextValue $ws "D40" "2.44"
Set-TextValue $ws "E40" "  +1.17%  "

# Row 41 - ARBITRUM
Set-TextValue $ws "E41" "  -0.53%  "

# Rows 42/43 swap rank: Aave now outranks WEMIXToken
Set-TextValue $ws "B42" "Aave"
Set-TextValue $ws "C42" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws "D42" "81.94"
Set-TextValue $ws "E42" "  +0.15%  "

Set-TextValue $ws "B43" "WEMIXToken"
Set-TextValue $ws "C43" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws "D43" "2.80"
Set-TextValue $ws "E43" "  -1.48%  "

# Row 44 - MXToken
Set-TextValue $ws "D44" "2.79"
Set-TextValue $ws "E44" "  -0.79%  "

# Row 45 - InjectiveProtocol
Set-TextValue $ws "E45" "  +1.62%  "

# Row 46 - Kaspa
Set-TextValue $ws "D46" "0.0509"
Set-TextValue $ws "E46" "  +1.54%  "

# Row 47 - RocketPoolETH
Set-TextValue $ws "D47" "1.974.13"
Set-TextValue $ws "E47" "  +0.48%  "

# Row 48 - FraxShare
Set-TextValue $ws "E48" "  -4.08%  "

# Row 49 - PaxDollar
Set-TextValue $ws "E49" "  -0.21%  "

# Row 50 - Quant
Set-TextValue $ws "D50" "102.65"
Set-TextValue $ws "E50" "  -2.22%  "

# Row 51 - BabyDogeCoin
Set-TextValue $ws "D51" "0.0₆0123"
Set-TextValue $ws "E51" "  -3.69%  "
